$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# setCellData fix: creating brand-new cells that did not previously exist
$ws.Range("E2").Value = "test"
$ws.Range("F6").Value = "test"
$ws.Range("F8").Value = "test2"

# D2 no longer carries the Hyperlink cell style
$ws.Range("D2").Style = "Normal"
$wb.Styles.Item("Hyperlink").Delete() | Out-Null

# selection moved to C5
$ws.Range("C5").Select() | Out-Null
